# Weekly update: insert a new "Mango" price-report row for Vega Modelo de
# Temuco ahead of the existing row 544, pushing all subsequent rows down by
# one (old row 544 becomes 545, ..., old row 614 becomes 615).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 544; Excel shifts rows 544:614 down to 545:615
# and the new row inherits the formatting (incl. the date style) of the row
# that used to sit there.
$ws.Rows.Item(544).Insert()

# Populate the newly inserted row 544 with the new record's data.
$ws.Range("A544").Value = 10
$ws.Range("B544").Value = "Vega Modelo de Temuco"
$ws.Range("C544").Value = "La Araucanía"
$ws.Range("D544").Value = 45131
$ws.Range("E544").Value = 9
$ws.Range("F544").Value = "Fruta"
$ws.Range("G544").Value = 100108
$ws.Range("H544").Value = "Tropicales y subtropicales"
$ws.Range("I544").Value = 100108002
$ws.Range("J544").Value = "Mango"
$ws.Range("K544").Value = "Sin especificar"
$ws.Range("L544").Value = "Primera"
$ws.Range("M544").Value = 1600
$ws.Range("N544").Value = 8000
$ws.Range("O544").Value = 9000
$ws.Range("P544").Value = 8562
$ws.Range("Q544").Value = "$/bandeja 4 kilos"
$ws.Range("R544").Value = "Brasil"
$ws.Range("S544").Value = 2140
$ws.Range("T544").Value = 4
